$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "AddCustomerTest"

$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

$ws.Range("A2").Value = "Katya"
$ws.Range("B2").Value = "Smith"
$ws.Range("C2").Value = "ab214c"

$ws.Range("C2").Select()
